# Adjustment Service Export Data Excel
#
# Restructures the "users" sheet from a simple A1:G7 Pre/Post-test summary
# table into a wider A1:AK9 layout with grouped Screening / Pretest /
# Posttest question headers (merged row-1 group labels + row-2 per-question
# labels) while keeping the original A:G summary table (now shifted down one
# row, with a new respondent appended at the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Push everything down one row so the old header row (A1:G1) becomes
#    A2:G2, the old "Vivi" row becomes row 3, ... and the old last row
#    ("ER") becomes row 8. Row 1 is now completely empty.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).Insert()

# ---------------------------------------------------------------------
# 2. Row 1 - merged group headers: Screening (H:I), Pretest (J:W),
#    Posttest (X:AK). A2 still carries the workbook's original bold
#    header style (it is the old header row, not cleared yet), so copy
#    its format into the whole H1:AK1 band before merging - this reuses
#    the existing bold style record instead of minting a new one.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy($ws.Range("H1:AK1")) | Out-Null

$ws.Range("H1").Value = "Screening"
$ws.Range("J1").Value = "Pretest"
$ws.Range("X1").Value = "Posttest"

$ws.Range("H1:I1").Merge()
$ws.Range("J1:W1").Merge()
$ws.Range("X1:AK1").Merge()

# ---------------------------------------------------------------------
# 3. Row 2 - per-question labels (plain/unstyled). Screening has 2
#    questions (H2:I2); Pretest and Posttest each repeat the same 14
#    question labels.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 8).Value = "Question 1"   # H2
$ws.Cells.Item(2, 9).Value = "Question 2"   # I2

$colIdx = 10  # J
foreach ($n in 1..14) {
    $ws.Cells.Item(2, $colIdx).Value = "Question $n"
    $colIdx++
}

$colIdx = 24  # X
foreach ($n in 1..14) {
    $ws.Cells.Item(2, $colIdx).Value = "Question $n"
    $colIdx++
}

# Now that the question labels are in place, wipe the stale A2:G2 cells
# left over from the old (now shifted) header row.
$ws.Range("A2:G2").Clear()

# ---------------------------------------------------------------------
# 4. Row 9 - append the new respondent ("Putria") after the existing
#    A:G summary table (rows 3-8, shifted down from the original rows
#    2-7, are already correct and untouched).
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value = "Putria"
$ws.Cells.Item(9, 2).Value = 22
$ws.Cells.Item(9, 3).Value = 23
$ws.Cells.Item(9, 4).Value = "S1"
$ws.Cells.Item(9, 5).Value = ">= 3 Juta"
$ws.Cells.Item(9, 6).Value = 18
$ws.Cells.Item(9, 7).Value = 15

Write-Output "done"
